$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from row 2 through row 496 to new date serial 45175
$ws.Range("C2:C496").Value = 45175
